# Applies the 8.2.1 "Annual growth rate of real GDP per employed person"
# update: adds a new 2022 column (S) and refreshes the 2020/2021 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header in S4, matching the format of the existing year headers
# (P4:R4, style index 21).
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("S4").Value = 2022

# Refresh existing data values for 2020 (Q5) and 2021 (R5), and add the new
# 2022 value (S5) using the same number format as R5 (style index 22).
$ws.Range("R5").Copy() | Out-Null
$ws.Range("Q5:S5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("Q5").Value = 91.892815141492093
$ws.Range("R5").Value = 101.53074848578628
$ws.Range("S5").Value = 109.27053140096621

$excel.CutCopyMode = 0

# Update the active selection like the recorded session did.
$ws.Range("T5").Select() | Out-Null
